# "mise à jour industrie"
# Fill in explicit zero values for previously-blank consumption/emission
# cells (columns D:J) in the Production_system sheet, and switch the
# active sheet/selection back to the "0D" tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Production_system")
$ws.Select()

# Row 2 (SMR): D,F,G,H,I were blank -> 0 (E already holds a formula)
$ws.Range("D2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0

# Row 3 (eSMR): F,G,H,I were blank -> 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0

# Row 4 (Electrolysis): E,F,G,H,I,J were blank -> 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0

# Row 5 (Coal-Gasification): E,F,G,H were blank -> 0 (I already a formula)
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0

# Row 6 (Biomass-Gasification): E,F,H,I were blank -> 0 (G already a formula)
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0

# Row 7 (Gas-Pyrolysis): F,G,H,I,J were blank -> 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0

# Update the saved selection on this sheet
$ws.Range("C14").Select()

# Make "0D" the active sheet/tab again
$ws0d = $wb.Worksheets.Item("0D")
$ws0d.Select()
